# Applies the cryptos.xlsx price/volume update described in the commit diff.
# Uses a leading apostrophe for values that look numeric, so Excel stores them
# as text (matching the original inlineStr/text cell type) instead of silently
# converting "0.630" -> 0.63, "1.00" -> 1, etc.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.913.30"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "2.241.64"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'246.25"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").Value = "'0.631"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "'74.59"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.630"
$ws.Range("E9").Value = "  -2.90%  "
$ws.Range("D10").Value = "'40.03"
$ws.Range("E10").Value = "  +0.96%  "
$ws.Range("D11").Value = "'0.0948"
$ws.Range("E11").Value = "  -4.34%  "
$ws.Range("D12").Value = "'7.17"
$ws.Range("E12").Value = "  -2.65%  "
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("D14").Value = "2.573.58"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").Value = "'14.83"
$ws.Range("E15").Value = "  -3.02%  "
$ws.Range("D16").Value = "'0.859"
$ws.Range("E16").Value = "  -2.46%  "
$ws.Range("D17").Value = "2.229.29"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").Value = "41.842.82"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").Value = "0.0₃0979"
$ws.Range("E19").Value = "  -2.74%  "
$ws.Range("D20").Value = "'6.12"
$ws.Range("E20").Value = "  -3.02%  "
$ws.Range("D21").Value = "'71.44"
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("D22").Value = "'2.24"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("D23").Value = "'230.50"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'11.25"
$ws.Range("E25").Value = "  -2.38%  "
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").Value = "'3.69"
$ws.Range("E26").Value = "  -5.90%  "
$ws.Range("D27").Value = "'2.31"
$ws.Range("E27").Value = "  -4.57%  "
$ws.Range("E28").Value = "  +12.60%  "
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D30").Value = "'168.84"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("D31").Value = "'20.52"
$ws.Range("E31").Value = "  -2.72%  "
$ws.Range("D32").Value = "'33.72"
$ws.Range("E32").Value = "  +3.73%  "
$ws.Range("D33").Value = "'0.0840"
$ws.Range("E33").Value = "  +2.69%  "
$ws.Range("E34").Value = "  -6.63%  "
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").Value = "'4.62"
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("D37").Value = "'4.89"
$ws.Range("E37").Value = "  +2.68%  "
$ws.Range("D38").Value = "'0.0299"
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("D39").Value = "'13.45"
$ws.Range("E39").Value = "  -4.53%  "
$ws.Range("D40").Value = "'5.91"
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("E41").Value = "  -6.79%  "
$ws.Range("D42").Value = "'110.85"
$ws.Range("E42").Value = "  +13.06%  "
$ws.Range("D43").Value = "'0.200"
$ws.Range("E43").Value = "  -6.59%  "
$ws.Range("D44").Value = "'60.24"
$ws.Range("E44").Value = "  -3.02%  "
$ws.Range("D45").Value = "'8.79"
$ws.Range("E45").Value = "  -3.74%  "
$ws.Range("D46").Value = "'0.101"
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("D47").Value = "'0.997"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("E48").Value = "  -4.29%  "
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("E50").Value = "  -12.74%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "'4.17"
$ws.Range("E51").Value = "  -2.89%  "
